# "Loan RBI, Variable Instalments"
#
# On the "Repayment schedule" sheet, insert a new (blank) column before the
# existing "Late" column (old column N). This pushes the old N ("Late"),
# O ("heading" spacer) and P ("Outstanding") columns one place to the right
# (-> O, P, Q) and widens the sheet's used range from A1:P8 to A1:Q8. The
# new column N is left blank (same cell style as the header/body cells it
# sits beside) and is given a fixed width of 11 characters (not auto
# "best fit", since it has no content to size itself to).
#
# Also mark "Repayment schedule" as the active/selected sheet (it was
# "Transactions" before), and leave the new cell selection at R10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column in front of column N ("Late"); shifts N,O,P -> O,P,Q
$ws.Columns("N").Insert() | Out-Null

# New column N has no best-fit content, so give it an explicit width that
# serializes to the workbook's stored column width of 11.
$ws.Columns("N").ColumnWidth = 10.166666666666666

# Make "Repayment schedule" the active sheet/tab, with R10 selected.
$ws.Activate()
$ws.Range("R10").Select() | Out-Null
